# Agregar Percepcion y deduccion
# Solo jala agregar por empleado, falta por departamento
#
# The "Fecha" column (column B) is no longer needed, so it is removed
# entirely - this shifts "Bono" (was C) into B and "Porcentaje" (was D)
# into C, and drops the date-format style that was only used by that
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Delete()

$ws.Range("F5").Select()
